$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 - this pushes the old row 10 ("Total")
# down to row 11, matching the diff's dimension change (A1:H10 -> A1:H11).
$ws.Rows.Item(10).Insert()

# Row 9 ("August (through 08-31)" -> "August") and its 2021 (column H) count.
$ws.Range("A9").Value = "August"
$ws.Range("H9").Value = 157

# New row 10: September (through 09-01) data.
$ws.Range("A10").Value = "September (through 09-01)"
# B10 stays blank for 2015 (no value yet for September), but still gets
# touched so a (empty) cell is materialized at B10, same as the source.
$ws.Range("B10").Font.Bold = $false
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 8
$ws.Range("G10").Value = 3
$ws.Range("H10").Value = 6

# The row-insert leaves A10 with a slightly different (border-less) style;
# restore the header-style formatting (bold/centered/bordered) used by the
# rest of column A by copying it from the neighboring row.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats

# Updated Total row, now row 11.
$ws.Range("B11").Value = 194
$ws.Range("C11").Value = 382
$ws.Range("D11").Value = 553
$ws.Range("E11").Value = 492
$ws.Range("F11").Value = 357
$ws.Range("G11").Value = 787
$ws.Range("H11").Value = 1076
